# Applies F-column numeric updates per the source diff, across three worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 307   # F3: 306 -> 307
$ws.Cells.Item(8, 6).Value = 1629   # F8: 1624 -> 1629
$ws.Cells.Item(10, 6).Value = 846   # F10: 845 -> 846
$ws.Cells.Item(11, 6).Value = 108   # F11: 106 -> 108
$ws.Cells.Item(14, 6).Value = 1496   # F14: 1493 -> 1496
$ws.Cells.Item(15, 6).Value = 6979   # F15: 6968 -> 6979
$ws.Cells.Item(17, 6).Value = 7144   # F17: 7140 -> 7144
$ws.Cells.Item(19, 6).Value = 5119   # F19: 4838 -> 5119
$ws.Cells.Item(20, 6).Value = 3091   # F20: 3089 -> 3091
$ws.Cells.Item(21, 6).Value = 3461   # F21: 3459 -> 3461
$ws.Cells.Item(23, 6).Value = 162   # F23: 158 -> 162
$ws.Cells.Item(24, 6).Value = 1842   # F24: 1841 -> 1842
$ws.Cells.Item(26, 6).Value = 293   # F26: 291 -> 293
$ws.Cells.Item(28, 6).Value = 18   # F28: 14 -> 18
$ws.Cells.Item(31, 6).Value = 2386   # F31: 1978 -> 2386
$ws.Cells.Item(32, 6).Value = 1124   # F32: 1119 -> 1124
$ws.Cells.Item(33, 6).Value = 2595   # F33: 2586 -> 2595
$ws.Cells.Item(34, 6).Value = 7   # F34: 5 -> 7
$ws.Cells.Item(36, 6).Value = 161   # F36: 160 -> 161
$ws.Cells.Item(37, 6).Value = 375   # F37: 374 -> 375
$ws.Cells.Item(38, 6).Value = 1030   # F38: 1028 -> 1030
$ws.Cells.Item(40, 6).Value = 467   # F40: 466 -> 467
$ws.Cells.Item(41, 6).Value = 515   # F41: 514 -> 515

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(8, 6).Value = 206   # F8: 205 -> 206

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 307   # F4: 306 -> 307
$ws.Cells.Item(9, 6).Value = 1629   # F9: 1624 -> 1629
$ws.Cells.Item(11, 6).Value = 846   # F11: 845 -> 846
$ws.Cells.Item(12, 6).Value = 108   # F12: 106 -> 108
$ws.Cells.Item(15, 6).Value = 1496   # F15: 1493 -> 1496
$ws.Cells.Item(16, 6).Value = 206   # F16: 205 -> 206
$ws.Cells.Item(19, 6).Value = 6980   # F19: 6968 -> 6980
$ws.Cells.Item(21, 6).Value = 7144   # F21: 7140 -> 7144
$ws.Cells.Item(23, 6).Value = 5119   # F23: 4839 -> 5119
$ws.Cells.Item(24, 6).Value = 3091   # F24: 3089 -> 3091
$ws.Cells.Item(25, 6).Value = 3461   # F25: 3459 -> 3461
$ws.Cells.Item(30, 6).Value = 1842   # F30: 1841 -> 1842
$ws.Cells.Item(33, 6).Value = 293   # F33: 291 -> 293
$ws.Cells.Item(35, 6).Value = 18   # F35: 14 -> 18
$ws.Cells.Item(38, 6).Value = 2386   # F38: 1979 -> 2386
$ws.Cells.Item(39, 6).Value = 1124   # F39: 1119 -> 1124
$ws.Cells.Item(41, 6).Value = 2595   # F41: 2586 -> 2595
$ws.Cells.Item(43, 6).Value = 161   # F43: 160 -> 161
$ws.Cells.Item(45, 6).Value = 375   # F45: 374 -> 375
$ws.Cells.Item(46, 6).Value = 1030   # F46: 1028 -> 1030
$ws.Cells.Item(48, 6).Value = 467   # F48: 466 -> 467
$ws.Cells.Item(49, 6).Value = 515   # F49: 514 -> 515
